$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C: plain text (names, URLs) - never ambiguous as numbers.
# Columns D/E: price/volume text that LOOKS numeric (e.g. "575.10", "2.521.39",
# "  -2.59%  "). The source workbook stores these as plain text, so each write
# is done as: mark the cell Text ("@"), assign the literal string, then clear
# the format again so the cell keeps style index 0 (matching the original) while
# the stored value remains the exact text (not auto-coerced into a float, which
# would silently drop things like a trailing zero or the 000-thousands dots).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.238.83'
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.75%  '
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.525.84'
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.94%  '
$ws.Range("E3").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.10'
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.80%  '
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.20'
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.71%  '
$ws.Range("E6").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.510'
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.75%  '
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.525.00'
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.00%  '
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.46%  '
$ws.Range("E10").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.343'
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.63%  '
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.81'
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.03%  '
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.990.18'
$ws.Range("D14").ClearFormats()

$ws.Range("B15").Value = 'ShibaInu'

$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000180'
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.41%  '
$ws.Range("E15").ClearFormats()

$ws.Range("B16").Value = 'WrappedBTC'

$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '70.056.15'
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.56%  '
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.90'
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.18%  '
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.532.23'
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -4.80%  '
$ws.Range("E18").ClearFormats()

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.43'
$ws.Range("D19").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -6.63%  '
$ws.Range("E19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("D20").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -8.25%  '
$ws.Range("E20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '353.59'
$ws.Range("D21").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.57%  '
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.93'
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -5.67%  '
$ws.Range("E22").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.29%  '
$ws.Range("E23").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '68.86'
$ws.Range("D25").ClearFormats()

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -4.34%  '
$ws.Range("E25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.07'
$ws.Range("D26").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -5.83%  '
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.23'
$ws.Range("D27").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.67%  '
$ws.Range("E27").ClearFormats()

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.655.43'
$ws.Range("D28").ClearFormats()

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.86%  '
$ws.Range("E28").ClearFormats()

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("E29").ClearFormats()

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0913'
$ws.Range("D30").ClearFormats()

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.81%  '
$ws.Range("E30").ClearFormats()

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.84'
$ws.Range("D31").ClearFormats()

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.93%  '
$ws.Range("E31").ClearFormats()

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '480.57'
$ws.Range("D32").ClearFormats()

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -4.12%  '
$ws.Range("E32").ClearFormats()

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.29'
$ws.Range("D33").ClearFormats()

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("E33").ClearFormats()

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.71%  '
$ws.Range("E34").ClearFormats()

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("E35").ClearFormats()

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.36'
$ws.Range("D36").ClearFormats()

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.53%  '
$ws.Range("E36").ClearFormats()

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.116'
$ws.Range("D37").ClearFormats()

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.12%  '
$ws.Range("E37").ClearFormats()

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("E38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.56'
$ws.Range("D39").ClearFormats()

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.79%  '
$ws.Range("E39").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("E40").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.67%  '
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.65'
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.82%  '
$ws.Range("E42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.319'
$ws.Range("D43").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.86%  '
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.72'
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -5.27%  '
$ws.Range("E44").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -7.05%  '
$ws.Range("E45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.31'
$ws.Range("D46").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.93%  '
$ws.Range("E46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.26'
$ws.Range("D47").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.84%  '
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.53'
$ws.Range("D48").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.47%  '
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.524'
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.25%  '
$ws.Range("E49").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -6.80%  '
$ws.Range("E50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.598'
$ws.Range("D51").ClearFormats()

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.17%  '
$ws.Range("E51").ClearFormats()
